$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell to a literal text value, forcing text format when the
# string would otherwise be auto-converted to a number by Excel (e.g. "316.33").
function Set-TextValue($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

$ws.Cells.Item(2, 4).Value = "42.096.23"
$ws.Cells.Item(2, 5).Value = "  -1.88%  "

$ws.Cells.Item(3, 4).Value = "2.292.41"
$ws.Cells.Item(3, 5).Value = "  -2.98%  "

$ws.Cells.Item(4, 5).Value = "  +0.06%  "

Set-TextValue 5 4 "316.33"
$ws.Cells.Item(5, 5).Value = "  +0.38%  "

Set-TextValue 6 4 "104.13"
$ws.Cells.Item(6, 5).Value = "  -4.08%  "

Set-TextValue 7 4 "0.630"

$ws.Cells.Item(8, 5).Value = "  +0.10%  "

$ws.Cells.Item(9, 5).Value = "  -2.41%  "

Set-TextValue 10 4 "39.53"
$ws.Cells.Item(10, 5).Value = "  -4.28%  "

$ws.Cells.Item(11, 5).Value = "  -2.27%  "

$ws.Cells.Item(12, 5).Value = "  -3.65%  "

$ws.Cells.Item(13, 5).Value = "  -0.27%  "

Set-TextValue 14 4 "0.961"
$ws.Cells.Item(14, 5).Value = "  -4.54%  "

Set-TextValue 15 4 "15.27"
$ws.Cells.Item(15, 5).Value = "  -4.48%  "

$ws.Cells.Item(16, 4).Value = "2.640.11"
$ws.Cells.Item(16, 5).Value = "  -3.14%  "

$ws.Cells.Item(17, 4).Value = "2.302.45"
$ws.Cells.Item(17, 5).Value = "  -2.41%  "

$ws.Cells.Item(18, 4).Value = "42.224.52"
$ws.Cells.Item(18, 5).Value = "  -1.57%  "

Set-TextValue 19 4 "7.34"
$ws.Cells.Item(19, 5).Value = "  -3.91%  "

Set-TextValue 20 4 "0.0000105"
$ws.Cells.Item(20, 5).Value = "  -0.99%  "

Set-TextValue 21 4 "73.33"
$ws.Cells.Item(21, 5).Value = "  -3.86%  "

$ws.Cells.Item(22, 5).Value = "  +0.19%  "

Set-TextValue 23 4 "277.88"
$ws.Cells.Item(23, 5).Value = "  +4.09%  "

$ws.Cells.Item(24, 5).Value = "  +10.21%  "

$ws.Cells.Item(25, 5).Value = "  -2.75%  "

$ws.Cells.Item(26, 5).Value = "  +0.79%  "

Set-TextValue 27 4 "10.82"
$ws.Cells.Item(27, 5).Value = "  -5.64%  "

$ws.Cells.Item(28, 5).Value = "  +5.05%  "

Set-TextValue 29 4 "22.78"
$ws.Cells.Item(29, 5).Value = "  -2.29%  "

Set-TextValue 30 4 "35.98"
$ws.Cells.Item(30, 5).Value = "  -2.46%  "

Set-TextValue 31 4 "163.78"
$ws.Cells.Item(31, 5).Value = "  -2.98%  "

Set-TextValue 32 4 "0.0872"
$ws.Cells.Item(32, 5).Value = "  -3.89%  "

Set-TextValue 33 4 "5.81"
$ws.Cells.Item(33, 5).Value = "  -4.82%  "

Set-TextValue 34 4 "2.81"
$ws.Cells.Item(34, 5).Value = "  -3.17%  "

$ws.Cells.Item(35, 5).Value = "  +3.53%  "

Set-TextValue 36 4 "0.112"
$ws.Cells.Item(36, 5).Value = "  -5.07%  "

Set-TextValue 37 4 "4.55"
$ws.Cells.Item(37, 5).Value = "  -3.31%  "

Set-TextValue 38 4 "0.0348"
$ws.Cells.Item(38, 5).Value = "  -3.95%  "

$ws.Cells.Item(39, 5).Value = "  -2.43%  "

Set-TextValue 40 4 "2.77"
$ws.Cells.Item(40, 5).Value = "  +3.35%  "

Set-TextValue 41 4 "99.85"
$ws.Cells.Item(41, 5).Value = "  -4.13%  "

$ws.Cells.Item(42, 5).Value = "  -4.13%  "

Set-TextValue 43 4 "69.22"
$ws.Cells.Item(43, 5).Value = "  -2.85%  "

# Row 44 <-> Row 45 swap: Algorand and FirstDigitalUSD exchange positions,
# and each also gets updated Price/Volume values.
$ws.Cells.Item(44, 2).Value = "FirstDigitalUSD"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue 44 4 "1.00"
$ws.Cells.Item(44, 5).Value = "  +0.18%  "

$ws.Cells.Item(45, 2).Value = "Algorand"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue 45 4 "0.225"
$ws.Cells.Item(45, 5).Value = "  -5.44%  "

Set-TextValue 46 4 "12.00"
$ws.Cells.Item(46, 5).Value = "  -4.79%  "

Set-TextValue 47 4 "111.88"
$ws.Cells.Item(47, 5).Value = "  -2.03%  "

Set-TextValue 48 4 "76.76"
$ws.Cells.Item(48, 5).Value = "  -5.75%  "

Set-TextValue 49 4 "8.92"
$ws.Cells.Item(49, 5).Value = "  -2.84%  "

Set-TextValue 50 4 "5.29"
$ws.Cells.Item(50, 5).Value = "  -5.12%  "

$ws.Cells.Item(51, 4).Value = "1.596.13"
$ws.Cells.Item(51, 5).Value = "  +0.79%  "
